$d = $word.ActiveDocument

function New-ParaXml($innerBody) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' + $innerBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1: change the "Wyplata dla czlonkow zespolu" paragraph: new text, drop the amount/tabs ---
$p4 = $d.Paragraphs(4)
$p4Xml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
  '<w:ind w:left="1440" w:hanging="359"/><w:contextualSpacing w:val="1"/><w:rPr/></w:pPr>' +
  '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
  '<w:t xml:space="preserve">Wypłata dla poszczególnych członków zespołu</w:t></w:r></w:p>'
$p4.Range.InsertXML((New-ParaXml $p4Xml))

# --- Step 2: insert the four new role/salary paragraphs right after it ---
$p4 = $d.Paragraphs(4)
$afterP4 = $p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5Xml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' +
  '<w:ind w:left="2160" w:hanging="359"/><w:contextualSpacing w:val="1"/><w:rPr/></w:pPr>' +
  '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
  '<w:t xml:space="preserve">Team leader of Software developers</w:t>' +
  '<w:tab/><w:tab/><w:tab/><w:tab/>' +
  '<w:t xml:space="preserve">9000</w:t><w:tab/><w:t xml:space="preserve">zł</w:t></w:r></w:p>'
$p5.Range.InsertXML((New-ParaXml $p5Xml))

$p5 = $d.Paragraphs(5)
$afterP5 = $p5.Range.InsertParagraphAfter()

$p6 = $d.Paragraphs(6)
$p6Xml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' +
  '<w:ind w:left="2160" w:hanging="359"/><w:contextualSpacing w:val="1"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
  '<w:t xml:space="preserve">Senior Software Engineer</w:t>' +
  '<w:tab/><w:tab/><w:tab/><w:tab/><w:tab/>' +
  '<w:t xml:space="preserve">7000</w:t><w:tab/><w:t xml:space="preserve">zł</w:t></w:r></w:p>'
$p6.Range.InsertXML((New-ParaXml $p6Xml))

$p6 = $d.Paragraphs(6)
$afterP6 = $p6.Range.InsertParagraphAfter()

$p7 = $d.Paragraphs(7)
$p7Xml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' +
  '<w:ind w:left="2160" w:hanging="359"/><w:contextualSpacing w:val="1"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
  '<w:t xml:space="preserve">Software Engineer</w:t>' +
  '<w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/>' +
  '<w:t xml:space="preserve">5000</w:t><w:tab/><w:t xml:space="preserve">zł</w:t></w:r></w:p>'
$p7.Range.InsertXML((New-ParaXml $p7Xml))

$p7 = $d.Paragraphs(7)
$afterP7 = $p7.Range.InsertParagraphAfter()

$p8 = $d.Paragraphs(8)
$p8Xml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' +
  '<w:ind w:left="2160" w:hanging="359"/><w:contextualSpacing w:val="1"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
  '<w:t xml:space="preserve">Junior Software Engineer</w:t>' +
  '<w:tab/><w:tab/><w:tab/><w:tab/><w:tab/>' +
  '<w:t xml:space="preserve">3000</w:t><w:tab/><w:t xml:space="preserve">zł</w:t></w:r></w:p>'
$p8.Range.InsertXML((New-ParaXml $p8Xml))

# --- Step 3: bump "Suma kosztow" 28400 -> 34400 (the cyclic-cost total, first occurrence) ---
# Paragraph index shifted by +4 (we inserted 4 new paragraphs above it).
$pSuma = $d.Paragraphs(11)
$sumaXml = '<w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:contextualSpacing w:val="0"/></w:pPr>' +
  '<w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr>' +
  '<w:t xml:space="preserve">Suma kosztów</w:t>' +
  '<w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/>' +
  '<w:t xml:space="preserve">34400</w:t><w:tab/><w:t xml:space="preserve">zł</w:t></w:r></w:p>'
$pSuma.Range.InsertXML((New-ParaXml $sumaXml))
